$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 20 continues the daily log with the next day after row 19
# (A19 = 12/8/2019 -> A20 = 12/9/2019, serial 43808).
# Copy the date cell's formatting from the row above so the new cell
# keeps the existing short-date style instead of minting a new number
# format.
$ws.Range("A19").Copy()
$ws.Range("A20").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("A20").Value = 43808
$ws.Range("B20:J20").Value = 0
$ws.Range("K20").Value = 1
$ws.Range("L20").Value = 2
$ws.Range("M20").Value = 2

# Leave the selection on K20, matching where data entry left off.
$ws.Range("K20").Select()
